# Apply the commit: update Fecha (D) and Volumen-derived (J,K,L,M,N,O,P,Q) columns
# for rows 48..165 by shifting each row's values down by one (row r <- old row r-1),
# inject new data at row 48 (D=45260, J=95), and append a new data row 166
# (carrying forward the values that shift out of row 165), matching the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    ,@(48,45260,95,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(49,45041,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(50,45051,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(51,44721,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(52,45159,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(53,45015,150,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(54,44679,90,12000,12000,12000,'$/caja 18 unidades','Región Metropolitana',667,18)
    ,@(55,45021,25,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(56,44385,100,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(57,45259,45,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(58,45054,200,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(59,45054,100,8000,8000,8000,'$/caja 18 unidades','Región del Maule',444,18)
    ,@(60,44778,45,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(61,44720,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(62,45061,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(63,45061,90,7000,7000,7000,'$/caja 18 unidades','Región del Maule',389,18)
    ,@(64,45076,85,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(65,45173,120,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(66,45204,65,12000,12000,12000,'$/caja 18 unidades','Región del Maule',667,18)
    ,@(67,45166,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(68,44741,180,10000,11000,10556,'$/caja 18 unidades','Región Metropolitana',586,18)
    ,@(69,44729,65,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
    ,@(70,45062,110,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(71,45085,200,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(72,45085,220,9000,10000,9455,'$/caja 18 unidades','Región del Maule',525,18)
    ,@(73,45127,200,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(74,45196,150,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(75,44706,150,11000,11000,11000,'$/caja 18 unidades','Región Metropolitana',611,18)
    ,@(76,45033,55,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(77,45075,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(78,45075,65,7500,7500,7500,'$/caja 18 unidades','Región del Maule',417,18)
    ,@(79,45167,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(80,45131,65,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(81,45246,55,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(82,45013,50,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(83,44386,40,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(84,44783,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(85,44775,40,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(86,45112,80,8500,8500,8500,'$/caja 18 unidades','Región Metropolitana',472,18)
    ,@(87,45106,150,8000,10000,9133,'$/caja 18 unidades','Región Metropolitana',507,18)
    ,@(88,44420,45,8000,8000,8000,'$/caja 16 unidades','Región Metropolitana',500,16)
    ,@(89,45043,250,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(90,45043,100,7500,8000,7750,'$/caja 18 unidades','Región del Maule',431,18)
    ,@(91,45030,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(92,44781,300,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(93,45146,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(94,45174,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(95,44771,100,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(96,45029,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(97,44762,85,11000,11000,11000,'$/caja 18 unidades','Región Metropolitana',611,18)
    ,@(98,45138,300,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(99,45162,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(100,45035,35,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(101,44392,95,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(102,44749,125,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
    ,@(103,45117,150,8000,8000,8000,'$/caja 18 unidades','Región del Maule',444,18)
    ,@(104,44348,35,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(105,45134,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(106,44847,110,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(107,45068,250,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(108,45068,200,8000,8000,8000,'$/caja 18 unidades','Región del Maule',444,18)
    ,@(109,45028,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(110,45169,300,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(111,45012,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(112,44719,50,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(113,45055,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(114,45069,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(115,44763,65,11000,11000,11000,'$/caja 18 unidades','Región Metropolitana',611,18)
    ,@(116,44756,550,10000,11000,10455,'$/caja 18 unidades','Región Metropolitana',581,18)
    ,@(117,44750,55,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
    ,@(118,45027,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(119,44354,100,8000,9000,8500,'$/caja 16 unidades','Región Metropolitana',531,16)
    ,@(120,44354,80,9000,9000,9000,'$/caja 16 unidades','Región del Maule',562,16)
    ,@(121,45225,150,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(122,44837,300,8000,8000,8000,'$/caja 18 unidades','Región del Maule',444,18)
    ,@(123,45163,90,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(124,45149,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(125,44369,60,7000,7000,7000,'$/caja 16 unidades','Región Metropolitana',438,16)
    ,@(126,44362,25,8000,8000,8000,'$/caja 16 unidades','Región Metropolitana',500,16)
    ,@(127,44315,40,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(128,44757,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(129,45141,300,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(130,44313,20,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(131,44767,500,10000,11000,10600,'$/caja 18 unidades','Región Metropolitana',589,18)
    ,@(132,45022,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(133,44816,65,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
    ,@(134,44740,170,10000,11000,10471,'$/caja 18 unidades','Región Metropolitana',582,18)
    ,@(135,44397,40,8000,8000,8000,'$/caja 16 unidades','Región Metropolitana',500,16)
    ,@(136,44764,45,11000,11000,11000,'$/caja 18 unidades','Región Metropolitana',611,18)
    ,@(137,45020,125,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(138,45244,35,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(139,45155,220,9000,10000,9455,'$/caja 18 unidades','Región Metropolitana',525,18)
    ,@(140,45040,100,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(141,45079,45,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(142,45096,55,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(143,45258,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(144,44676,40,12000,12000,12000,'$/caja 18 unidades','Región Metropolitana',667,18)
    ,@(145,45034,40,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(146,44848,45,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(147,45243,110,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(148,45044,80,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(149,45126,100,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(150,45007,25,10000,10000,10000,'$/caja 18 unidades','Región del Maule',556,18)
    ,@(151,45099,110,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(152,44685,90,12000,12000,12000,'$/caja 18 unidades','Región del Maule',667,18)
    ,@(153,44312,40,7000,7000,7000,'$/caja 16 unidades','Región del Maule',438,16)
    ,@(154,44777,65,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(155,44396,80,7000,7000,7000,'$/caja 16 unidades','Región Metropolitana',438,16)
    ,@(156,45222,400,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(157,44761,55,11000,11000,11000,'$/caja 18 unidades','Región Metropolitana',611,18)
    ,@(158,45097,55,9000,9000,9000,'$/caja 18 unidades','Región Metropolitana',500,18)
    ,@(159,44399,80,7000,7000,7000,'$/caja 16 unidades','Región Metropolitana',438,16)
    ,@(160,45050,35,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(161,45050,65,8000,8000,8000,'$/caja 18 unidades','Región del Maule',444,18)
    ,@(162,45180,150,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(163,44726,125,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
    ,@(164,44727,35,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(165,45135,75,10000,10000,10000,'$/caja 18 unidades','Región Metropolitana',556,18)
    ,@(166,45100,65,8000,8000,8000,'$/caja 18 unidades','Región Metropolitana',444,18)
)

# Row 166 needs the same constant columns (A,B,C,E,F,G,H,I,R) as every other data row.
# Copy them from row 165 before row 165's own D..Q values get overwritten.
$ws.Cells.Item(166,1).Value = $ws.Cells.Item(165,1).Value2
$ws.Cells.Item(166,2).Value = $ws.Cells.Item(165,2).Value2
$ws.Cells.Item(166,3).Value = $ws.Cells.Item(165,3).Value2
$ws.Cells.Item(166,5).Value = $ws.Cells.Item(165,5).Value2
$ws.Cells.Item(166,6).Value = $ws.Cells.Item(165,6).Value2
$ws.Cells.Item(166,7).Value = $ws.Cells.Item(165,7).Value2
$ws.Cells.Item(166,8).Value = $ws.Cells.Item(165,8).Value2
$ws.Cells.Item(166,9).Value = $ws.Cells.Item(165,9).Value2
$ws.Cells.Item(166,18).Value = $ws.Cells.Item(165,18).Value2

# New row's date cell (D166) needs the same date/time number format as the rest of column D.
$ws.Cells.Item(166,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value  = $entry[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $entry[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $entry[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $entry[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $entry[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $entry[6]   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $entry[7]   # O - Origen
    $ws.Cells.Item($r, 16).Value = $entry[8]   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $entry[9]   # Q - Kg o Unidades
}
